# Applies the cryptos list refresh described in the commit
# "Updated cryptos list on Fri Sep 15 05:49:01 UTC 2023 with GitHub Actions".
# Only the Price (D) and Volume(1h) (E) columns change; everything else is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price strings that are plain decimals (e.g. "213.61") would otherwise be
# auto-coerced into numbers by Excel's type inference when assigned via
# .Value, losing the original text formatting (and picking up float noise,
# e.g. 213.61000000000001). Forcing the cell to text format first keeps them
# as literal strings, matching the source data (multi-dot prices such as
# "26.674.20" are never parsed as numbers, so they don't need this).
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

$ws.Range("D2").Value = '26.674.20'
$ws.Range("D3").Value = '1.635.38'
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").Value = '  +0.04%  '
Set-TextValue $ws.Range("D5") '213.61'
$ws.Range("E5").Value = '  +0.79%  '
Set-TextValue $ws.Range("D6") '0.508'
$ws.Range("E6").Value = '  +4.43%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +2.60%  '
$ws.Range("E9").Value = '  +1.53%  '
Set-TextValue $ws.Range("D10") '19.26'
$ws.Range("E10").Value = '  +2.86%  '
$ws.Range("E11").Value = '  +3.68%  '
$ws.Range("D12").Value = '1.866.19'
$ws.Range("E12").Value = '  +1.23%  '
$ws.Range("D13").Value = '1.620.57'
$ws.Range("E13").Value = '  +0.06%  '
Set-TextValue $ws.Range("D14") '4.10'
$ws.Range("E14").Value = '  +2.75%  '
$ws.Range("E15").Value = '  +2.10%  '
$ws.Range("D16").Value = '26.686.42'
$ws.Range("E16").Value = '  +1.57%  '
Set-TextValue $ws.Range("D17") '63.55'
$ws.Range("E17").Value = '  +2.16%  '
$ws.Range("D18").Value = '0.0₃0745'
$ws.Range("E18").Value = '  +2.59%  '
Set-TextValue $ws.Range("D19") '219.58'
$ws.Range("E19").Value = '  +9.43%  '
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("E21").Value = '  +1.03%  '
Set-TextValue $ws.Range("D22") '9.45'
$ws.Range("E22").Value = '  +1.63%  '
Set-TextValue $ws.Range("D23") '6.20'
$ws.Range("E23").Value = '  +2.80%  '
Set-TextValue $ws.Range("D24") '1.92'
$ws.Range("E24").Value = '  +1.35%  '
Set-TextValue $ws.Range("D25") '148.16'
$ws.Range("E25").Value = '  +2.81%  '
$ws.Range("E26").Value = '  +0.00%  '
Set-TextValue $ws.Range("D27") '0.121'
$ws.Range("E27").Value = '  +1.66%  '
Set-TextValue $ws.Range("D28") '6.93'
$ws.Range("E28").Value = '  +6.20%  '
Set-TextValue $ws.Range("D29") '15.49'
$ws.Range("E29").Value = '  +2.26%  '
Set-TextValue $ws.Range("D30") '0.0511'
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("E31").Value = '  +0.13%  '
Set-TextValue $ws.Range("D32") '3.32'
$ws.Range("E32").Value = '  +4.95%  '
$ws.Range("E33").Value = '  +2.71%  '
$ws.Range("E34").Value = '  +1.94%  '
$ws.Range("E35").Value = '  -0.31%  '
$ws.Range("D36").Value = '1.214.02'
$ws.Range("E36").Value = '  +3.32%  '
$ws.Range("E37").Value = '  +5.81%  '
Set-TextValue $ws.Range("D38") '0.812'
$ws.Range("E38").Value = '  +1.38%  '
Set-TextValue $ws.Range("D40") '0.507'
$ws.Range("E40").Value = '  +2.83%  '
$ws.Range("E41").Value = '  -1.10%  '
Set-TextValue $ws.Range("D42") '5.44'
$ws.Range("E42").Value = '  +2.25%  '
Set-TextValue $ws.Range("D43") '0.794'
$ws.Range("E43").Value = '  +0.26%  '
$ws.Range("D44").Value = '1.775.71'
$ws.Range("E44").Value = '  +1.19%  '
Set-TextValue $ws.Range("D45") '93.34'
$ws.Range("E45").Value = '  +0.84%  '
Set-TextValue $ws.Range("D46") '1.55'
$ws.Range("E46").Value = '  +2.29%  '
Set-TextValue $ws.Range("D47") '54.92'
$ws.Range("E47").Value = '  +2.53%  '
$ws.Range("E48").Value = '  +1.05%  '
Set-TextValue $ws.Range("D49") '7.69'
$ws.Range("E49").Value = '  +6.08%  '
$ws.Range("E50").Value = '  +0.42%  '
$ws.Range("E51").Value = '  +0.07%  '
